# Updated cryptos list on Thu Mar 28 17:39:59 UTC 2024 with GitHub Actions
# Applies refreshed price/volume data (and a Stellar/THORChain row swap)
# to the cryptos worksheet. Values are written as plain text (not numbers),
# matching the original inlineStr cell encoding exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.737.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.553.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.533.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E10").Value = "  +17.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.652"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000315"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.117.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.698.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.574.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "570.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.91%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -14.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.28%  "
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "546.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.415"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.75%  "
$ws.Range("E38").Value = "  +9.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.573.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.42%  "
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0451"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.138"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
